# Updated cryptos list on Tue Aug 20 18:42:03 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose values would otherwise be
# auto-converted to numbers by Excel (values look numeric, but source
# data is stored as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.205.15"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "2.583.15"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "568.19"
$ws.Range("E5").Value = "  +2.43%  "

$ws.Range("D6").Value = "143.28"
$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").Value = "2.588.00"
$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("E11").Value = "  +2.92%  "

$ws.Range("E12").Value = "  +8.51%  "

$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").Value = "3.037.93"
$ws.Range("E14").Value = "  -0.46%  "

$ws.Range("D15").Value = "59.243.30"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("D16").Value = "22.49"
$ws.Range("E16").Value = "  +7.67%  "

$ws.Range("E17").Value = "  +4.08%  "

$ws.Range("D18").Value = "2.589.09"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").Value = "4.51"
$ws.Range("E19").Value = "  +1.60%  "

$ws.Range("D20").Value = "335.87"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").Value = "10.20"
$ws.Range("E21").Value = "  +1.36%  "

$ws.Range("D22").Value = "6.18"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").Value = "64.07"
$ws.Range("E24").Value = "  -3.48%  "

$ws.Range("D25").Value = "0.452"
$ws.Range("E25").Value = "  +5.53%  "

$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  +1.46%  "

$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  +3.63%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "1.68"
$ws.Range("E31").Value = "  +0.31%  "

$ws.Range("D32").Value = "6.06"

$ws.Range("D33").Value = "157.34"
$ws.Range("E33").Value = "  +2.78%  "

$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("D35").Value = "4.03"
$ws.Range("E35").Value = "  +2.78%  "

$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").Value = "0.880"
$ws.Range("E36").Value = "  +0.88%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "0.885"
$ws.Range("E37").Value = "  +7.33%  "

$ws.Range("E38").Value = "  +1.93%  "

$ws.Range("E39").Value = "  +2.90%  "

$ws.Range("D40").Value = "36.86"
$ws.Range("E40").Value = "  -0.33%  "

$ws.Range("D41").Value = "294.31"
$ws.Range("E41").Value = "  +3.74%  "

$ws.Range("D42").Value = "3.66"
$ws.Range("E42").Value = "  +1.21%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "0.0973"
$ws.Range("E44").Value = "  +1.89%  "

$ws.Range("D45").Value = "0.597"
$ws.Range("E45").Value = "  -0.26%  "

$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "10.63"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "19.15"
$ws.Range("E48").Value = "  +2.39%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "124.49"
$ws.Range("E49").Value = "  +5.27%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0232"
$ws.Range("E50").Value = "  +2.48%  "

$ws.Range("D51").Value = "18.50"
$ws.Range("E51").Value = "  +3.86%  "

# Restore default (General) formatting/style on the price cells now that
# their text values are set, without altering any other formatting.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
